$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, $value)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "244.35"
Set-TextValue $ws.Range("D3") "21.87"
Set-TextValue $ws.Range("D4") "5.394"
Set-TextValue $ws.Range("D5") "0.05997"
Set-TextValue $ws.Range("D8") "0.9500"
Set-TextValue $ws.Range("D9") "0.0005902"
Set-TextValue $ws.Range("E9") "8OneONE"
Set-TextValue $ws.Range("D10") "0.1427"
Set-TextValue $ws.Range("D11") "0.07397"
Set-TextValue $ws.Range("D12") "0.03310"
Set-TextValue $ws.Range("D14") "0.09406"
Set-TextValue $ws.Range("D15") "4.003"
Set-TextValue $ws.Range("D16") "0.001594"
Set-TextValue $ws.Range("D17") "0.04820"
Set-TextValue $ws.Range("D18") "0.006222"
Set-TextValue $ws.Range("D19") "0.005002"
Set-TextValue $ws.Range("D20") "0.0009917"
Set-TextValue $ws.Range("D23") "6.415"
Set-TextValue $ws.Range("D25") "0.3253"
Set-TextValue $ws.Range("D26") "0.1340"
Set-TextValue $ws.Range("D40") "0.03990"
Set-TextValue $ws.Range("D41") "0.006513"
Set-TextValue $ws.Range("D42") "0.1072"
Set-TextValue $ws.Range("D44") "0.005246"
Set-TextValue $ws.Range("D45") "0.00005257"
Set-TextValue $ws.Range("E47") "46CoinbaseStockTokenCOINBestin24h"
Set-TextValue $ws.Range("D48") "0.01620"
